$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: 10uF Ceramic Capacitor -> add ", 25V" and fill in LCSC part number
$ws.Range("A3").Value = "10uF Ceramic Capacitor, 25V"
$ws.Range("D3").Value = "C96446"

# Row 4: 100nF Ceramic Capacitor - fill in LCSC part number (reuses C1591)
$ws.Range("D4").Value = "C1591"

# Row 5: 22uF Ceramic Capacitor - fill in LCSC part number
$ws.Range("D5").Value = "C2762594"

# Row 6: 3.3uH Inductor - fill in Footprint
$ws.Range("C6").Value = "SMD,7.2x6.6mm"

# Row 9: 56.2k Resistor -> add ", 0.1%" and fill in LCSC part number
$ws.Range("A9").Value = "56.2k Resistor, 0.1%"
$ws.Range("D9").Value = "C705784"

# Row 10: 10k Resistor -> add ", 0.1%" and fill in LCSC part number
$ws.Range("A10").Value = "10k Resistor, 0.1%"
$ws.Range("D10").Value = "C95204"

# Update the saved selection to match the author's final cursor position
$ws.Range("A10").Select()
